$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "ECs" sending-cluster row); rows 3 and 4 shift up to become rows 2 and 3.
$ws.Rows(2).Delete()

# After the delete, the former row 3 (FAPs) is now row 2, and the former row 4 (MuSCs) is now row 3.
# Update the recalculated specificity columns (I, J, S, T) for the new row 2 and row 3.
$ws.Range("I2").Value = 0.5234889777227022
$ws.Range("J2").Value = 0.5234889777227022
$ws.Range("S2").Value = 0.5234889777227022
$ws.Range("T2").Value = 0.5234889777227022

$ws.Range("I3").Value = 0.4765110222772978
$ws.Range("J3").Value = 0.4765110222772979
$ws.Range("S3").Value = 0.4765110222772978
$ws.Range("T3").Value = 0.4765110222772979
